$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 7).Value2 = 3.263122
$ws.Cells.Item(2, 8).Value2 = 9.789365999999999
$ws.Cells.Item(2, 9).Value2 = 0.3531375780718168
$ws.Cells.Item(2, 10).Value2 = 0.3531375780718168
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 13).Value2 = 127.6999736666667
$ws.Cells.Item(2, 14).Value2 = 383.099921
$ws.Cells.Item(2, 15).Value2 = 0.9554352891750322
$ws.Cells.Item(2, 16).Value2 = 0.9554352891750322
$ws.Cells.Item(2, 17).Value2 = 416.7005934711206
$ws.Cells.Item(2, 18).Value2 = 3750.305341240086
$ws.Cells.Item(2, 19).Value2 = 0.3374001040236168
$ws.Cells.Item(2, 20).Value2 = 0.3374001040236168

$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 7).Value2 = 3.263122
$ws.Cells.Item(3, 8).Value2 = 9.789365999999999
$ws.Cells.Item(3, 9).Value2 = 0.3531375780718168
$ws.Cells.Item(3, 10).Value2 = 0.3531375780718168
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 13).Value2 = 0.4321196666666667
$ws.Cells.Item(3, 14).Value2 = 1.296359
$ws.Cells.Item(3, 15).Value2 = 0.003233065495828321
$ws.Cells.Item(3, 16).Value2 = 0.003233065495828321
$ws.Cells.Item(3, 17).Value2 = 1.410059190932667
$ws.Cells.Item(3, 18).Value2 = 12.690532718394
$ws.Cells.Item(3, 19).Value2 = 0.001141716918944371
$ws.Cells.Item(3, 20).Value2 = 0.001141716918944371

$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 7).Value2 = 3.263122
$ws.Cells.Item(4, 8).Value2 = 9.789365999999999
$ws.Cells.Item(4, 9).Value2 = 0.3531375780718168
$ws.Cells.Item(4, 10).Value2 = 0.3531375780718168
$ws.Cells.Item(4, 11).Value2 = 3
$ws.Cells.Item(4, 13).Value2 = 4.77305
$ws.Cells.Item(4, 14).Value2 = 14.31915
$ws.Cells.Item(4, 15).Value2 = 0.03571136528892854
$ws.Cells.Item(4, 16).Value2 = 0.03571136528892854
$ws.Cells.Item(4, 17).Value2 = 15.5750444621
$ws.Cells.Item(4, 18).Value2 = 140.1754001589
$ws.Cells.Item(4, 19).Value2 = 0.01261102504777017
$ws.Cells.Item(4, 20).Value2 = 0.01261102504777017

$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 7).Value2 = 3.263122
$ws.Cells.Item(5, 8).Value2 = 9.789365999999999
$ws.Cells.Item(5, 9).Value2 = 0.3531375780718168
$ws.Cells.Item(5, 10).Value2 = 0.3531375780718168
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 13).Value2 = 0.751186
$ws.Cells.Item(5, 14).Value2 = 2.253558
$ws.Cells.Item(5, 15).Value2 = 0.00562028004021099
$ws.Cells.Item(5, 16).Value2 = 0.00562028004021099
$ws.Cells.Item(5, 17).Value2 = 2.451211562692
$ws.Cells.Item(5, 18).Value2 = 22.060904064228
$ws.Cells.Item(5, 19).Value2 = 0.001984732081485482
$ws.Cells.Item(5, 20).Value2 = 0.001984732081485482

$ws.Cells.Item(6, 5).Value2 = 3
$ws.Cells.Item(6, 7).Value2 = 4.367310666666667
$ws.Cells.Item(6, 8).Value2 = 13.101932
$ws.Cells.Item(6, 9).Value2 = 0.4726337266929886
$ws.Cells.Item(6, 10).Value2 = 0.4726337266929886
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 13).Value2 = 127.6999736666667
$ws.Cells.Item(6, 14).Value2 = 383.099921
$ws.Cells.Item(6, 15).Value2 = 0.9554352891750322
$ws.Cells.Item(6, 16).Value2 = 0.9554352891750322
$ws.Cells.Item(6, 17).Value2 = 557.7054571274857
$ws.Cells.Item(6, 18).Value2 = 5019.349114147371
$ws.Cells.Item(6, 19).Value2 = 0.4515709413367887
$ws.Cells.Item(6, 20).Value2 = 0.4515709413367887

$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 7).Value2 = 4.367310666666667
$ws.Cells.Item(7, 8).Value2 = 13.101932
$ws.Cells.Item(7, 9).Value2 = 0.4726337266929886
$ws.Cells.Item(7, 10).Value2 = 0.4726337266929886
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 13).Value2 = 0.4321196666666667
$ws.Cells.Item(7, 14).Value2 = 1.296359
$ws.Cells.Item(7, 15).Value2 = 0.003233065495828321
$ws.Cells.Item(7, 16).Value2 = 0.003233065495828321
$ws.Cells.Item(7, 17).Value2 = 1.887200829509778
$ws.Cells.Item(7, 18).Value2 = 16.984807465588
$ws.Cells.Item(7, 19).Value2 = 0.001528055793935854
$ws.Cells.Item(7, 20).Value2 = 0.001528055793935854

$ws.Cells.Item(8, 5).Value2 = 3
$ws.Cells.Item(8, 7).Value2 = 4.367310666666667
$ws.Cells.Item(8, 8).Value2 = 13.101932
$ws.Cells.Item(8, 9).Value2 = 0.4726337266929886
$ws.Cells.Item(8, 10).Value2 = 0.4726337266929886
$ws.Cells.Item(8, 11).Value2 = 3
$ws.Cells.Item(8, 13).Value2 = 4.77305
$ws.Cells.Item(8, 14).Value2 = 14.31915
$ws.Cells.Item(8, 15).Value2 = 0.03571136528892854
$ws.Cells.Item(8, 16).Value2 = 0.03571136528892854
$ws.Cells.Item(8, 17).Value2 = 20.84539217753333
$ws.Cells.Item(8, 18).Value2 = 187.6085295978
$ws.Cells.Item(8, 19).Value2 = 0.01687839566180093
$ws.Cells.Item(8, 20).Value2 = 0.01687839566180093

$ws.Cells.Item(9, 5).Value2 = 3
$ws.Cells.Item(9, 7).Value2 = 4.367310666666667
$ws.Cells.Item(9, 8).Value2 = 13.101932
$ws.Cells.Item(9, 9).Value2 = 0.4726337266929886
$ws.Cells.Item(9, 10).Value2 = 0.4726337266929886
$ws.Cells.Item(9, 11).Value2 = 3
$ws.Cells.Item(9, 13).Value2 = 0.751186
$ws.Cells.Item(9, 14).Value2 = 2.253558
$ws.Cells.Item(9, 15).Value2 = 0.00562028004021099
$ws.Cells.Item(9, 16).Value2 = 0.00562028004021099
$ws.Cells.Item(9, 17).Value2 = 3.280662630450667
$ws.Cells.Item(9, 18).Value2 = 29.525963674056
$ws.Cells.Item(9, 19).Value2 = 0.00265633390046314
$ws.Cells.Item(9, 20).Value2 = 0.00265633390046314

$ws.Cells.Item(10, 5).Value2 = 3
$ws.Cells.Item(10, 7).Value2 = 1.609937666666666
$ws.Cells.Item(10, 8).Value2 = 4.829813
$ws.Cells.Item(10, 9).Value2 = 0.1742286952351946
$ws.Cells.Item(10, 10).Value2 = 0.1742286952351946
$ws.Cells.Item(10, 11).Value2 = 3
$ws.Cells.Item(10, 13).Value2 = 127.6999736666667
$ws.Cells.Item(10, 14).Value2 = 383.099921
$ws.Cells.Item(10, 15).Value2 = 0.9554352891750322
$ws.Cells.Item(10, 16).Value2 = 0.9554352891750322
$ws.Cells.Item(10, 17).Value2 = 205.5889976383081
$ws.Cells.Item(10, 18).Value2 = 1850.300978744773
$ws.Cells.Item(10, 19).Value2 = 0.1664642438146267
$ws.Cells.Item(10, 20).Value2 = 0.1664642438146267

$ws.Cells.Item(11, 5).Value2 = 3
$ws.Cells.Item(11, 7).Value2 = 1.609937666666666
$ws.Cells.Item(11, 8).Value2 = 4.829813
$ws.Cells.Item(11, 9).Value2 = 0.1742286952351946
$ws.Cells.Item(11, 10).Value2 = 0.1742286952351946
$ws.Cells.Item(11, 11).Value2 = 3
$ws.Cells.Item(11, 13).Value2 = 0.4321196666666667
$ws.Cells.Item(11, 14).Value2 = 1.296359
$ws.Cells.Item(11, 15).Value2 = 0.003233065495828321
$ws.Cells.Item(11, 16).Value2 = 0.003233065495828321
$ws.Cells.Item(11, 17).Value2 = 0.6956857278741111
$ws.Cells.Item(11, 18).Value2 = 6.261171550867
$ws.Cells.Item(11, 19).Value2 = 0.0005632927829480957
$ws.Cells.Item(11, 20).Value2 = 0.0005632927829480957

$ws.Cells.Item(12, 5).Value2 = 3
$ws.Cells.Item(12, 7).Value2 = 1.609937666666666
$ws.Cells.Item(12, 8).Value2 = 4.829813
$ws.Cells.Item(12, 9).Value2 = 0.1742286952351946
$ws.Cells.Item(12, 10).Value2 = 0.1742286952351946
$ws.Cells.Item(12, 11).Value2 = 3
$ws.Cells.Item(12, 13).Value2 = 4.77305
$ws.Cells.Item(12, 14).Value2 = 14.31915
$ws.Cells.Item(12, 15).Value2 = 0.03571136528892854
$ws.Cells.Item(12, 16).Value2 = 0.03571136528892854
$ws.Cells.Item(12, 17).Value2 = 7.684312979883333
$ws.Cells.Item(12, 18).Value2 = 69.15881681895
$ws.Cells.Item(12, 19).Value2 = 0.006221944579357436
$ws.Cells.Item(12, 20).Value2 = 0.006221944579357436

$ws.Cells.Item(13, 5).Value2 = 3
$ws.Cells.Item(13, 7).Value2 = 1.609937666666666
$ws.Cells.Item(13, 8).Value2 = 4.829813
$ws.Cells.Item(13, 9).Value2 = 0.1742286952351946
$ws.Cells.Item(13, 10).Value2 = 0.1742286952351946
$ws.Cells.Item(13, 11).Value2 = 3
$ws.Cells.Item(13, 13).Value2 = 0.751186
$ws.Cells.Item(13, 14).Value2 = 2.253558
$ws.Cells.Item(13, 15).Value2 = 0.00562028004021099
$ws.Cells.Item(13, 16).Value2 = 0.00562028004021099
$ws.Cells.Item(13, 17).Value2 = 1.209362636072667
$ws.Cells.Item(13, 18).Value2 = 10.884263724654
$ws.Cells.Item(13, 19).Value2 = 0.0009792140582623678
$ws.Cells.Item(13, 20).Value2 = 0.0009792140582623678
